$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing log entry: "opeteltuauth" -> "opeteltu auth"
$ws.Cells.Item(66, 3).Value = "opeteltu auth  cookien lähettäminen suoraan serveriltä, testattu ja apin:n koodi refaktoroitu"

# New row 67: new time-log entry for the "api" project
$ws.Cells.Item(67, 2).Value = 3
$ws.Cells.Item(67, 3).Value = "api herokuun, github action automaattiseen liven päivittämiseen, production versiossa oma mongodb url(ei omaa konttia)"
$ws.Cells.Item(67, 4).Value = "api"

# Extend the totals formula to include the new row
$ws.Cells.Item(75, 2).Formula = "=SUM(B2:B67)"

# Update view state to match where the user ended up working
$ws.Range("B66").Select()

$excel.CalculateFull()
